$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "preview" block (before the old
# row 4 "dataset.commit.id"), pushing everything else down by two rows.
$ws.Rows("4:5").Insert()

# New key/value pairs describing the dataset preview queries.
$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"

$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

# The two new rows hold multi-line formulas, so make them tall and wrap
# their text.
$ws.Rows("4:5").RowHeight = 120
$ws.Range("A4:B5").VerticalAlignment = -4108
$ws.Range("A4:B5").WrapText = $true

# Match the author's final selection.
$ws.Range("B11").Select() | Out-Null
